$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player, position, team) replacing the old A2:C18 block.
$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The old table had one extra row (18); remove it entirely so the used range shrinks.
$newLastRow = 1 + $data.Count
$ws.Rows.Item($newLastRow + 1).Delete()
